# Updates the cryptos worksheet cell values to match the refreshed snapshot
# (coin prices / 1h volume %, plus a handful of rows whose coin identity and
# link/price/volume were swapped with the adjoining row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign as literal text. If the text looks like a plain number (e.g. "5.43"),
# prefix it with a leading apostrophe so Excel stores it as text instead of
# auto-converting it to a numeric value, matching the original inlineStr cells.
$numericPattern = '^[+-]?[0-9]*\.?[0-9]+$'

$val = '55.645.04'
if ($val -match $numericPattern) { $ws.Range("D2").Value = "'" + $val } else { $ws.Range("D2").Value = $val }
$val = '  +2.68%  '
if ($val -match $numericPattern) { $ws.Range("E2").Value = "'" + $val } else { $ws.Range("E2").Value = $val }
$val = '2.493.80'
if ($val -match $numericPattern) { $ws.Range("D3").Value = "'" + $val } else { $ws.Range("D3").Value = $val }
$val = '  +6.55%  '
if ($val -match $numericPattern) { $ws.Range("E3").Value = "'" + $val } else { $ws.Range("E3").Value = $val }
$val = '  +0.18%  '
if ($val -match $numericPattern) { $ws.Range("E4").Value = "'" + $val } else { $ws.Range("E4").Value = $val }
$val = '480.01'
if ($val -match $numericPattern) { $ws.Range("D5").Value = "'" + $val } else { $ws.Range("D5").Value = $val }
$val = '  +7.07%  '
if ($val -match $numericPattern) { $ws.Range("E5").Value = "'" + $val } else { $ws.Range("E5").Value = $val }
$val = '139.11'
if ($val -match $numericPattern) { $ws.Range("D6").Value = "'" + $val } else { $ws.Range("D6").Value = $val }
$val = '  +8.32%  '
if ($val -match $numericPattern) { $ws.Range("E6").Value = "'" + $val } else { $ws.Range("E6").Value = $val }
$val = '0.999'
if ($val -match $numericPattern) { $ws.Range("D7").Value = "'" + $val } else { $ws.Range("D7").Value = $val }
$val = '  +0.47%  '
if ($val -match $numericPattern) { $ws.Range("E7").Value = "'" + $val } else { $ws.Range("E7").Value = $val }
$val = '0.512'
if ($val -match $numericPattern) { $ws.Range("D8").Value = "'" + $val } else { $ws.Range("D8").Value = $val }
$val = '  +7.19%  '
if ($val -match $numericPattern) { $ws.Range("E8").Value = "'" + $val } else { $ws.Range("E8").Value = $val }
$val = '2.484.42'
if ($val -match $numericPattern) { $ws.Range("D9").Value = "'" + $val } else { $ws.Range("D9").Value = $val }
$val = '  +6.57%  '
if ($val -match $numericPattern) { $ws.Range("E9").Value = "'" + $val } else { $ws.Range("E9").Value = $val }
$val = '0.0981'
if ($val -match $numericPattern) { $ws.Range("D10").Value = "'" + $val } else { $ws.Range("D10").Value = $val }
$val = '  +5.73%  '
if ($val -match $numericPattern) { $ws.Range("E10").Value = "'" + $val } else { $ws.Range("E10").Value = $val }
$val = '5.43'
if ($val -match $numericPattern) { $ws.Range("D11").Value = "'" + $val } else { $ws.Range("D11").Value = $val }
$val = '  +0.69%  '
if ($val -match $numericPattern) { $ws.Range("E11").Value = "'" + $val } else { $ws.Range("E11").Value = $val }
$val = '0.326'
if ($val -match $numericPattern) { $ws.Range("D12").Value = "'" + $val } else { $ws.Range("D12").Value = $val }
$val = '  +4.91%  '
if ($val -match $numericPattern) { $ws.Range("E12").Value = "'" + $val } else { $ws.Range("E12").Value = $val }
$val = '  +0.47%  '
if ($val -match $numericPattern) { $ws.Range("E13").Value = "'" + $val } else { $ws.Range("E13").Value = $val }
$val = '2.923.27'
if ($val -match $numericPattern) { $ws.Range("D14").Value = "'" + $val } else { $ws.Range("D14").Value = $val }
$val = '  +7.97%  '
if ($val -match $numericPattern) { $ws.Range("E14").Value = "'" + $val } else { $ws.Range("E14").Value = $val }
$val = '55.617.80'
if ($val -match $numericPattern) { $ws.Range("D15").Value = "'" + $val } else { $ws.Range("D15").Value = $val }
$val = '  +2.68%  '
if ($val -match $numericPattern) { $ws.Range("E15").Value = "'" + $val } else { $ws.Range("E15").Value = $val }
$val = 'ShibaInu'
if ($val -match $numericPattern) { $ws.Range("B16").Value = "'" + $val } else { $ws.Range("B16").Value = $val }
$val = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
if ($val -match $numericPattern) { $ws.Range("C16").Value = "'" + $val } else { $ws.Range("C16").Value = $val }
$val = '0.0000137'
if ($val -match $numericPattern) { $ws.Range("D16").Value = "'" + $val } else { $ws.Range("D16").Value = $val }
$val = '  +12.29%  '
if ($val -match $numericPattern) { $ws.Range("E16").Value = "'" + $val } else { $ws.Range("E16").Value = $val }
$val = 'Avalanche'
if ($val -match $numericPattern) { $ws.Range("B17").Value = "'" + $val } else { $ws.Range("B17").Value = $val }
$val = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
if ($val -match $numericPattern) { $ws.Range("C17").Value = "'" + $val } else { $ws.Range("C17").Value = $val }
$val = '20.41'
if ($val -match $numericPattern) { $ws.Range("D17").Value = "'" + $val } else { $ws.Range("D17").Value = $val }
$val = '  +7.40%  '
if ($val -match $numericPattern) { $ws.Range("E17").Value = "'" + $val } else { $ws.Range("E17").Value = $val }
$val = '2.489.32'
if ($val -match $numericPattern) { $ws.Range("D18").Value = "'" + $val } else { $ws.Range("D18").Value = $val }
$val = '  +6.90%  '
if ($val -match $numericPattern) { $ws.Range("E18").Value = "'" + $val } else { $ws.Range("E18").Value = $val }
$val = '4.36'
if ($val -match $numericPattern) { $ws.Range("D19").Value = "'" + $val } else { $ws.Range("D19").Value = $val }
$val = '  +8.52%  '
if ($val -match $numericPattern) { $ws.Range("E19").Value = "'" + $val } else { $ws.Range("E19").Value = $val }
$val = '320.06'
if ($val -match $numericPattern) { $ws.Range("D20").Value = "'" + $val } else { $ws.Range("D20").Value = $val }
$val = '  +5.91%  '
if ($val -match $numericPattern) { $ws.Range("E20").Value = "'" + $val } else { $ws.Range("E20").Value = $val }
$val = '9.99'
if ($val -match $numericPattern) { $ws.Range("D21").Value = "'" + $val } else { $ws.Range("D21").Value = $val }
$val = '  +6.59%  '
if ($val -match $numericPattern) { $ws.Range("E21").Value = "'" + $val } else { $ws.Range("E21").Value = $val }
$val = '0.998'
if ($val -match $numericPattern) { $ws.Range("D22").Value = "'" + $val } else { $ws.Range("D22").Value = $val }
$val = '  -0.05%  '
if ($val -match $numericPattern) { $ws.Range("E22").Value = "'" + $val } else { $ws.Range("E22").Value = $val }
$val = '  +4.44%  '
if ($val -match $numericPattern) { $ws.Range("E23").Value = "'" + $val } else { $ws.Range("E23").Value = $val }
$val = '57.79'
if ($val -match $numericPattern) { $ws.Range("D24").Value = "'" + $val } else { $ws.Range("D24").Value = $val }
$val = '  +3.71%  '
if ($val -match $numericPattern) { $ws.Range("E24").Value = "'" + $val } else { $ws.Range("E24").Value = $val }
$val = '  +0.35%  '
if ($val -match $numericPattern) { $ws.Range("E25").Value = "'" + $val } else { $ws.Range("E25").Value = $val }
$val = '0.404'
if ($val -match $numericPattern) { $ws.Range("D26").Value = "'" + $val } else { $ws.Range("D26").Value = $val }
$val = '  +8.83%  '
if ($val -match $numericPattern) { $ws.Range("E26").Value = "'" + $val } else { $ws.Range("E26").Value = $val }
$val = '0.163'
if ($val -match $numericPattern) { $ws.Range("D27").Value = "'" + $val } else { $ws.Range("D27").Value = $val }
$val = '  +6.29%  '
if ($val -match $numericPattern) { $ws.Range("E27").Value = "'" + $val } else { $ws.Range("E27").Value = $val }
$val = '2.603.55'
if ($val -match $numericPattern) { $ws.Range("D28").Value = "'" + $val } else { $ws.Range("D28").Value = $val }
$val = '  +7.31%  '
if ($val -match $numericPattern) { $ws.Range("E28").Value = "'" + $val } else { $ws.Range("E28").Value = $val }
$val = '7.37'
if ($val -match $numericPattern) { $ws.Range("D29").Value = "'" + $val } else { $ws.Range("D29").Value = $val }
$val = '  +6.08%  '
if ($val -match $numericPattern) { $ws.Range("E29").Value = "'" + $val } else { $ws.Range("E29").Value = $val }
$val = '  +7.79%  '
if ($val -match $numericPattern) { $ws.Range("E30").Value = "'" + $val } else { $ws.Range("E30").Value = $val }
$val = '  +0.50%  '
if ($val -match $numericPattern) { $ws.Range("E31").Value = "'" + $val } else { $ws.Range("E31").Value = $val }
$val = '148.76'
if ($val -match $numericPattern) { $ws.Range("D32").Value = "'" + $val } else { $ws.Range("D32").Value = $val }
$val = '  +1.18%  '
if ($val -match $numericPattern) { $ws.Range("E32").Value = "'" + $val } else { $ws.Range("E32").Value = $val }
$val = '18.11'
if ($val -match $numericPattern) { $ws.Range("D33").Value = "'" + $val } else { $ws.Range("D33").Value = $val }
$val = '  +5.72%  '
if ($val -match $numericPattern) { $ws.Range("E33").Value = "'" + $val } else { $ws.Range("E33").Value = $val }
$val = '  +8.15%  '
if ($val -match $numericPattern) { $ws.Range("E34").Value = "'" + $val } else { $ws.Range("E34").Value = $val }
$val = '  +10.52%  '
if ($val -match $numericPattern) { $ws.Range("E35").Value = "'" + $val } else { $ws.Range("E35").Value = $val }
$val = '  +1.67%  '
if ($val -match $numericPattern) { $ws.Range("E36").Value = "'" + $val } else { $ws.Range("E36").Value = $val }
$val = '  +8.68%  '
if ($val -match $numericPattern) { $ws.Range("E37").Value = "'" + $val } else { $ws.Range("E37").Value = $val }
$val = '0.843'
if ($val -match $numericPattern) { $ws.Range("D38").Value = "'" + $val } else { $ws.Range("D38").Value = $val }
$val = '  -0.06%  '
if ($val -match $numericPattern) { $ws.Range("E38").Value = "'" + $val } else { $ws.Range("E38").Value = $val }
$val = '34.17'
if ($val -match $numericPattern) { $ws.Range("D39").Value = "'" + $val } else { $ws.Range("D39").Value = $val }
$val = '  +2.05%  '
if ($val -match $numericPattern) { $ws.Range("E39").Value = "'" + $val } else { $ws.Range("E39").Value = $val }
$val = '0.997'
if ($val -match $numericPattern) { $ws.Range("D40").Value = "'" + $val } else { $ws.Range("D40").Value = $val }
$val = '  +0.57%  '
if ($val -match $numericPattern) { $ws.Range("E40").Value = "'" + $val } else { $ws.Range("E40").Value = $val }
$val = '0.608'
if ($val -match $numericPattern) { $ws.Range("D41").Value = "'" + $val } else { $ws.Range("D41").Value = $val }
$val = '  +15.91%  '
if ($val -match $numericPattern) { $ws.Range("E41").Value = "'" + $val } else { $ws.Range("E41").Value = $val }
$val = '0.0549'
if ($val -match $numericPattern) { $ws.Range("D42").Value = "'" + $val } else { $ws.Range("D42").Value = $val }
$val = '  +10.10%  '
if ($val -match $numericPattern) { $ws.Range("E42").Value = "'" + $val } else { $ws.Range("E42").Value = $val }
$val = 'Stacks'
if ($val -match $numericPattern) { $ws.Range("B43").Value = "'" + $val } else { $ws.Range("B43").Value = $val }
$val = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
if ($val -match $numericPattern) { $ws.Range("C43").Value = "'" + $val } else { $ws.Range("C43").Value = $val }
$val = '1.32'
if ($val -match $numericPattern) { $ws.Range("D43").Value = "'" + $val } else { $ws.Range("D43").Value = $val }
$val = '  +6.82%  '
if ($val -match $numericPattern) { $ws.Range("E43").Value = "'" + $val } else { $ws.Range("E43").Value = $val }
$val = 'Filecoin'
if ($val -match $numericPattern) { $ws.Range("B44").Value = "'" + $val } else { $ws.Range("B44").Value = $val }
$val = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
if ($val -match $numericPattern) { $ws.Range("C44").Value = "'" + $val } else { $ws.Range("C44").Value = $val }
$val = '3.37'
if ($val -match $numericPattern) { $ws.Range("D44").Value = "'" + $val } else { $ws.Range("D44").Value = $val }
$val = '  +6.41%  '
if ($val -match $numericPattern) { $ws.Range("E44").Value = "'" + $val } else { $ws.Range("E44").Value = $val }
$val = '10.15'
if ($val -match $numericPattern) { $ws.Range("D45").Value = "'" + $val } else { $ws.Range("D45").Value = $val }
$val = '  -1.10%  '
if ($val -match $numericPattern) { $ws.Range("E45").Value = "'" + $val } else { $ws.Range("E45").Value = $val }
$val = 'Stellar'
if ($val -match $numericPattern) { $ws.Range("B46").Value = "'" + $val } else { $ws.Range("B46").Value = $val }
$val = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
if ($val -match $numericPattern) { $ws.Range("C46").Value = "'" + $val } else { $ws.Range("C46").Value = $val }
$val = '0.0906'
if ($val -match $numericPattern) { $ws.Range("D46").Value = "'" + $val } else { $ws.Range("D46").Value = $val }
$val = '  +9.11%  '
if ($val -match $numericPattern) { $ws.Range("E46").Value = "'" + $val } else { $ws.Range("E46").Value = $val }
$val = 'Maker'
if ($val -match $numericPattern) { $ws.Range("B47").Value = "'" + $val } else { $ws.Range("B47").Value = $val }
$val = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
if ($val -match $numericPattern) { $ws.Range("C47").Value = "'" + $val } else { $ws.Range("C47").Value = $val }
$val = '1.961.17'
if ($val -match $numericPattern) { $ws.Range("D47").Value = "'" + $val } else { $ws.Range("D47").Value = $val }
$val = '  +0.88%  '
if ($val -match $numericPattern) { $ws.Range("E47").Value = "'" + $val } else { $ws.Range("E47").Value = $val }
$val = 'VeChain'
if ($val -match $numericPattern) { $ws.Range("B48").Value = "'" + $val } else { $ws.Range("B48").Value = $val }
$val = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
if ($val -match $numericPattern) { $ws.Range("C48").Value = "'" + $val } else { $ws.Range("C48").Value = $val }
$val = '0.0223'
if ($val -match $numericPattern) { $ws.Range("D48").Value = "'" + $val } else { $ws.Range("D48").Value = $val }
$val = '  +5.56%  '
if ($val -match $numericPattern) { $ws.Range("E48").Value = "'" + $val } else { $ws.Range("E48").Value = $val }
$val = 'Bittensor'
if ($val -match $numericPattern) { $ws.Range("B49").Value = "'" + $val } else { $ws.Range("B49").Value = $val }
$val = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
if ($val -match $numericPattern) { $ws.Range("C49").Value = "'" + $val } else { $ws.Range("C49").Value = $val }
$val = '248.14'
if ($val -match $numericPattern) { $ws.Range("D49").Value = "'" + $val } else { $ws.Range("D49").Value = $val }
$val = '  +29.82%  '
if ($val -match $numericPattern) { $ws.Range("E49").Value = "'" + $val } else { $ws.Range("E49").Value = $val }
$val = '4.45'
if ($val -match $numericPattern) { $ws.Range("D50").Value = "'" + $val } else { $ws.Range("D50").Value = $val }
$val = '  +9.64%  '
if ($val -match $numericPattern) { $ws.Range("E50").Value = "'" + $val } else { $ws.Range("E50").Value = $val }
$val = '17.44'
if ($val -match $numericPattern) { $ws.Range("D51").Value = "'" + $val } else { $ws.Range("D51").Value = $val }
$val = '  +8.27%  '
if ($val -match $numericPattern) { $ws.Range("E51").Value = "'" + $val } else { $ws.Range("E51").Value = $val }
